# Weekly update: insert a new price record as row 329 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 329-411 down to 330-412.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 329; Excel copies the
# formatting (incl. the date number format in column D) from the row above.
$ws.Rows.Item(329).Insert()

# Populate the new row with the latest market reading.
$ws.Cells.Item(329, 1).Value = 4
$ws.Cells.Item(329, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(329, 3).Value = "Los Lagos"
$ws.Cells.Item(329, 4).Value = 44932
$ws.Cells.Item(329, 5).Value = 10
$ws.Cells.Item(329, 6).Value = 100114014
$ws.Cells.Item(329, 7).Value = "Betarraga"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 1000
$ws.Cells.Item(329, 11).Value = 1000
$ws.Cells.Item(329, 12).Value = 1000
$ws.Cells.Item(329, 13).Value = 1000
$ws.Cells.Item(329, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(329, 15).Value = "Región del Maule"
$ws.Cells.Item(329, 16).Value = 200
$ws.Cells.Item(329, 17).Value = 5
$ws.Cells.Item(329, 18).Value = "Hortaliza"
